$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.244.85'
$ws.Range('E2').Value = '  +1.08%  '

$ws.Range('D3').Value = '2.647.86'
$ws.Range('E3').Value = '  +3.34%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = '593.45'
$ws.Range('E5').Value = '  +2.67%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.30'
$ws.Range('E6').Value = '  +0.08%  '

$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').Value = '0.587'
$ws.Range('E8').Value = '  +0.12%  '

$ws.Range('D9').Value = '2.646.66'
$ws.Range('E9').Value = '  +3.38%  '

$ws.Range('E10').Value = '  +0.53%  '

$ws.Range('D11').Value = '5.67'
$ws.Range('E11').Value = '  +2.47%  '

$ws.Range('E12').Value = '  +0.83%  '

$ws.Range('E13').Value = '  +1.83%  '

$ws.Range('D14').Value = '27.37'
$ws.Range('E14').Value = '  +2.29%  '

$ws.Range('D15').Value = '3.126.58'
$ws.Range('E15').Value = '  +3.51%  '

$ws.Range('D16').Value = '63.212.45'
$ws.Range('E16').Value = '  +1.15%  '

$ws.Range('E17').Value = '  +0.46%  '

$ws.Range('D18').Value = '2.641.50'
$ws.Range('E18').Value = '  +3.57%  '

$ws.Range('D19').Value = '11.41'
$ws.Range('E19').Value = '  +3.00%  '

$ws.Range('D20').Value = '339.38'
$ws.Range('E20').Value = '  +0.56%  '

$ws.Range('D21').Value = '4.37'
$ws.Range('E21').Value = '  +1.29%  '

$ws.Range('D22').Value = '6.77'
$ws.Range('E22').Value = '  +2.11%  '

$ws.Range('E23').Value = '  +0.10%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '67.30'
$ws.Range('E24').Value = '  +0.31%  '

$ws.Range('E25').Value = '  +6.77%  '

$ws.Range('D26').Value = '0.165'
$ws.Range('E26').Value = '  +1.59%  '

$ws.Range('E27').Value = '  +1.29%  '

$ws.Range('E28').Value = '  +0.18%  '

$ws.Range('D29').Value = '8.44'
$ws.Range('E29').Value = '  +3.45%  '

$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').Value = '534.78'
$ws.Range('E30').Value = '  +17.39%  '

$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').Value = '7.83'
$ws.Range('E31').Value = '  -0.98%  '

$ws.Range('E32').Value = '  +14.20%  '

$ws.Range('E33').Value = '  +4.07%  '

$ws.Range('D34').Value = '0.0₃0807'
$ws.Range('E34').Value = '  +1.97%  '

$ws.Range('D35').Value = '173.95'
$ws.Range('E35').Value = '  -1.76%  '

$ws.Range('D36').Value = '5.09'
$ws.Range('E36').Value = '  +14.91%  '

$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').Value = '0.406'
$ws.Range('E37').Value = '  +3.07%  '

$ws.Range('B38').Value = 'FirstDigitalUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.01%  '

$ws.Range('D39').Value = '19.03'
$ws.Range('E39').Value = '  +1.36%  '

$ws.Range('D40').Value = '1.82'
$ws.Range('E40').Value = '  +8.55%  '

$ws.Range('D41').Value = '172.49'
$ws.Range('E41').Value = '  +9.84%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'

$ws.Range('D43').Value = '40.08'
$ws.Range('E43').Value = '  -1.00%  '

$ws.Range('E44').Value = '  +2.36%  '

$ws.Range('D45').Value = '22.06'
$ws.Range('E45').Value = '  +6.42%  '

$ws.Range('D46').Value = '0.0559'
$ws.Range('E46').Value = '  +5.37%  '

$ws.Range('D47').Value = '0.632'
$ws.Range('E47').Value = '  +0.88%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0240'
$ws.Range('E48').Value = '  +2.98%  '

$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0960'
$ws.Range('E49').Value = '  +0.58%  '

$ws.Range('D50').Value = '18.67'
$ws.Range('E50').Value = '  +4.48%  '

$ws.Range('E51').Value = '  +2.94%  '
